$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'331.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.20%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'41.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'3.89%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.685"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.61%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.08362"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.18%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'8.804"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.79%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'2.011"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.68%"
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'1.00%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'2.989"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.22%"
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'-0.31%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.1291"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.68%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.1969"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.93%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.09470"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'3.26%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.03881"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'8.68%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.1060"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.87%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.001306"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.64%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'0.006104"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-3.16%"
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'3.439"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.90%"
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.3539"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'2.63%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'8.016"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-7.98%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'0.1371"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.20%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.2609"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.15%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04413"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.05%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001255"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.26%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.004450"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.14%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0001202"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-2.70%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D39").Value = "'0.02801"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'3.98%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.05560"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.14%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.007976"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'6.13%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.1436"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.43%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.009298"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-5.78%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.002160"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'2.85%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.01110"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-4.60%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00007009"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'3.55%"
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'0.19%"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.003534"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'15.47%"
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'0.15%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.00002100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.19%"
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.0002000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.19%"
$ws.Range("E51").Style = "Normal"
